# "Update automatico via Actualizar 02-12-2021 15-51-25"
#
# This workbook keeps a rolling availability log: every refresh cycle
# shifts the previous two blocks of 14 rows' timestamps (column D) one
# block down, and stamps the first block (rows 2-15) with the current
# refresh time. Columns A-C (name/link/status) are untouched because the
# same 14-service sequence repeats in every block.
#
# Oldest block (rows 30-43) is overwritten with what used to be in
# rows 16-29, the middle block (rows 16-29) is overwritten with what
# used to be in rows 2-15, and the newest block (rows 2-15) gets the
# fresh "now" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Oldest block <- previous middle block's old timestamp
$ws.Range("D30:D43").Value2 = 44234.72344277778

# Middle block <- previous newest block's old timestamp
$ws.Range("D16:D29").Value2 = 44239.63938445602

# Newest block <- freshly captured update timestamp
$ws.Range("D2:D15").Value2 = 44239.66063602377
